$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in shared string used by A23 label
$ws.Range("A23").Value = "actual actuation force penalty per patch (N)"

# Update height of patch (mm) from 25 to 50
$ws.Range("B5").Value = 50

# Add new inlay-depth verification formula in J23
$ws.Range("J23").Formula = "=70.65+0.65*2+3.175"

# Update sheet view to match new location/selection
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("M13").Select()
